# Update the multiplication equations in the document to match the new
# "commit" (c8c62b6) output. Each old equation text is unique within the
# document, so a simple Find/Replace per pair is safe and unambiguous.

$d = $word.ActiveDocument

$replacements = @(
    @("16×78=1248", "47×17=799"),
    @("30×59=1770", "50×26=1300"),
    @("74×27=1998", "87×31=2697"),
    @("43×80=3440", "39×90=3510"),
    @("72×38=2736", "17×35=595"),
    @("87×47=4089", "54×57=3078"),
    @("17×84=1428", "85×66=5610"),
    @("80×84=6720", "59×71=4189"),
    @("55×28=1540", "73×50=3650"),
    @("86×87=7482", "45×24=1080"),
    @("99×62=6138", "90×69=6210"),
    @("38×93=3534", "29×44=1276"),
    @("20×68=1360", "95×97=9215"),
    @("54×79=4266", "42×53=2226"),
    @("29×16=464",  "47×94=4418"),
    @("85×62=5270", "19×40=760"),
    @("25×36=900",  "65×45=2925"),
    @("22×23=506",  "58×23=1334"),
    @("46×76=3496", "45×47=2115"),
    @("56×68=3808", "50×44=2200"),
    @("31×36=1116", "67×73=4891"),
    @("69×80=5520", "86×28=2408"),
    @("26×40=1040", "28×55=1540"),
    @("77×57=4389", "83×51=4233"),
    @("54×54=2916", "65×60=3900")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}
